# Update results for other subsets (dev other / test other WER columns)
# across the per-configuration result sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "general" ----
$ws = $wb.Worksheets.Item("general")
$ws.Range("E4").Value = 3
$ws.Range("H4").Value = 7.7
$ws.Range("E5").Value = 2.9
$ws.Range("G5").Value = 3.3
$ws.Range("H5").Value = 7.2
$ws.Range("E6").Value = 2.9
$ws.Range("H6").Value = 7.6
$ws.Range("E7").Value = 2.9
$ws.Range("H7").Value = 7.5
$ws.Range("E7").Select()

# ---- Sheet "pretraining" ----
$ws = $wb.Worksheets.Item("pretraining")
$ws.Range("C4").Value = 2.9
$ws.Range("F4").Value = 7.5
$ws.Range("C5").Value = 2.9
$ws.Range("F5").Value = 7.5
$ws.Range("C6").Value = 3
$ws.Range("F6").Value = 7.8
$ws.Range("F7").Select()

# ---- Sheet "window_size" ----
$ws = $wb.Worksheets.Item("window_size")
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 2.9
$ws.Range("F6").Value = 7.6
$ws.Range("C7").Value = 3
$ws.Range("F7").Value = 7.7
$ws.Range("C8").Value = 2.9
$ws.Range("C5").Select()

# ---- Sheet "scf_size" ----
$ws = $wb.Worksheets.Item("scf_size")
$ws.Range("E4").Value = 3
$ws.Range("H4").Value = 7.8
$ws.Range("E5").Value = 3
$ws.Range("H5").Value = 7.7
$ws.Range("E6").Value = 2.9
$ws.Range("E7").Value = 3
$ws.Range("H7").Value = 7.8
$ws.Range("E8").Value = 3
$ws.Range("H8").Value = 7.7
$ws.Range("E9").Value = 3
$ws.Range("H9").Value = 7.7
$ws.Range("H10").Select()

# ---- Sheet "w2v_size" ----
$ws = $wb.Worksheets.Item("w2v_size")
$ws.Range("D4").Value = 2.9
$ws.Range("G4").Value = 7.5
$ws.Range("D5").Value = 2.9
$ws.Range("G5").Value = 7.5
$ws.Range("D6").Value = 2.9
$ws.Range("G6").Value = 7.6
$ws.Range("D7").Value = 3
$ws.Range("G7").Value = 7.7
$ws.Range("D8").Value = 3.1
$ws.Range("G8").Value = 7.6
$ws.Range("D9").Value = 2.9
$ws.Range("D10").Value = 3
$ws.Range("G10").Value = 7.9
$ws.Range("D11").Value = 2.9
$ws.Range("G11").Value = 7.6
$ws.Range("D12").Value = 3.1
$ws.Range("D13").Value = 2.9
$ws.Range("G13").Value = 7.4
$ws.Range("D14").Value = 3.1
$ws.Range("G14").Value = 7.9
$ws.Range("D15").Value = 3
$ws.Range("G15").Value = 7.7
$ws.Range("D16").Value = 3.1
$ws.Range("G16").Select()

# ---- Sheet "w2v_proj" ----
$ws = $wb.Worksheets.Item("w2v_proj")
$ws.Range("C4").Value = 2.9
$ws.Range("C5").Value = 2.9
$ws.Range("F5").Value = 7.5
$ws.Range("F6").Select()

# Restore the originally active sheet/tab ("w2v_size" was tab-selected,
# activeTab index 4) after touching the other sheets above.
$wsActive = $wb.Worksheets.Item("w2v_size")
$wsActive.Activate()
$wsActive.Range("G16").Select()
